$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: rename the bean from "Effect" to "EventEffect", set separator, and
# change the first field name/type to "EffectType"
$ws.Range("B4").Value = "EventEffect"
$ws.Range("E4").Value = ","
$ws.Range("H4").Value = "EffectType"
$ws.Range("I4").Value = "EffectType"

# Row 5: new field Para1 : int
$ws.Range("H5").Value = "Para1"
$ws.Range("I5").Value = "int"

# Row 6: new field Para2 : int
$ws.Range("H6").Value = "Para2"
$ws.Range("I6").Value = "int"

# Touch B9 so it carries the "bad" style like its neighbours, matching the
# row-block pattern used elsewhere in the sheet. Copy formats only (no
# value) from B4, which already carries that style.
$ws.Range("B4").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to I7, mirroring the author's cursor position
# when the edit was saved.
$ws.Range("I7").Select()
